# Filter testcases for sales order added
#
# This script mutates the single worksheet ("Input_Data") of the workbook to
# match the target revision:
#   - Row 3 ("sales_orders" test case) gets its IMOrderNo / OrderType /
#     OrderStatus values corrected, plus TotalRevenueMin/TotalRevenueMax and
#     three brand-new filter columns (ResellerName / EndUserName / CreatedOn).
#   - Three new header columns are appended (V/W/X) with the matching header
#     style copied from the existing header row.
#   - A brand-new row 9 stashes the previous IMOrderNo value that used to
#     live in C3.
#   - The saved window/selection state is refreshed (F11 instead of U11,
#     no frozen/scrolled topLeftCell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New header cells: V1/W1/X1 (copy the header style from S1, then set text)
# ---------------------------------------------------------------------------
$ws.Range("V1").Value = "ResellerName"
$ws.Range("S1").Copy()
$ws.Range("V1").PasteSpecial(-4122)

$ws.Range("W1").Value = "EndUserName"
$ws.Range("S1").Copy()
$ws.Range("W1").PasteSpecial(-4122)

$ws.Range("X1").Value = "CreatedOn"
$ws.Range("S1").Copy()
$ws.Range("X1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# New column widths for the two freshly-introduced columns (V, W)
# ---------------------------------------------------------------------------
$ws.Columns.Item(22).ColumnWidth = 12.6
$ws.Columns.Item(23).ColumnWidth = 14.6

# ---------------------------------------------------------------------------
# Row 3 ("sales_orders" sample row) value corrections
# ---------------------------------------------------------------------------
# The old IMOrderNo value moves down to the new row 9 ...
$ws.Range("C9").Value = "20-VN1CR-11"
# ... and C3 takes the new order number used by the sales-order testcase.
$ws.Range("C3").Value = "20-VN2W9-11"

# OrderType corrected from the literal "S" to the full "Stock" value.
$ws.Range("D3").Value = "Stock"

# OrderStatus corrected from "Order Hold(IM)" to "Order Hold".
$ws.Range("G3").Value = "Order Hold"

# New TotalRevenueMin / TotalRevenueMax values for the sales-order filter.
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 300

# New ResellerName / EndUserName / CreatedOn filter values.
$ws.Range("V3").Value = "INGRAM MICRO CAP TEST ACCOUNT"
$ws.Range("W3").Value = "Everest EndUser"
$ws.Range("X3").Value = "Yesterday"

# ---------------------------------------------------------------------------
# Refresh the saved view/selection state
# ---------------------------------------------------------------------------
$ws.Range("F11").Select()
